$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Content.Find.Execute("2025-01-03 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-01-04 Saturday", 2)

# Update the five data rows (1, 5, 9, 13, 17) of the single table, 5 columns each.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "82÷2=41, 0"
$t.Cell(1, 2).Range.Text = "19÷2=9, 1"
$t.Cell(1, 3).Range.Text = "23÷7=3, 2"
$t.Cell(1, 4).Range.Text = "54÷9=6, 0"
$t.Cell(1, 5).Range.Text = "84÷6=14, 0"

$t.Cell(5, 1).Range.Text = "43÷6=7, 1"
$t.Cell(5, 2).Range.Text = "21÷2=10, 1"
$t.Cell(5, 3).Range.Text = "29÷2=14, 1"
$t.Cell(5, 4).Range.Text = "96÷8=12, 0"
$t.Cell(5, 5).Range.Text = "85÷6=14, 1"

$t.Cell(9, 1).Range.Text = "39÷6=6, 3"
$t.Cell(9, 2).Range.Text = "92÷7=13, 1"
$t.Cell(9, 3).Range.Text = "90÷7=12, 6"
$t.Cell(9, 4).Range.Text = "79÷9=8, 7"
$t.Cell(9, 5).Range.Text = "64÷7=9, 1"

$t.Cell(13, 1).Range.Text = "15÷7=2, 1"
$t.Cell(13, 2).Range.Text = "79÷7=11, 2"
$t.Cell(13, 3).Range.Text = "50÷2=25, 0"
$t.Cell(13, 4).Range.Text = "64÷8=8, 0"
$t.Cell(13, 5).Range.Text = "37÷4=9, 1"

$t.Cell(17, 1).Range.Text = "24÷5=4, 4"
$t.Cell(17, 2).Range.Text = "31÷6=5, 1"
$t.Cell(17, 3).Range.Text = "97÷2=48, 1"
$t.Cell(17, 4).Range.Text = "76÷5=15, 1"
$t.Cell(17, 5).Range.Text = "61÷9=6, 7"
